$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 53
$ws.Range("C2").Value = 55
$ws.Range("E2").Value = 49.07407407407408
$ws.Range("F2").Value = 0.330027
$ws.Range("G2").Value = 0.010541
$ws.Range("H2").Value = 0.001447917704583131
$ws.Range("I2").Value = 0.002837918700982937
$ws.Range("J2").Value = 0.3328649187009829
$ws.Range("K2").Value = 0.3271890812990171
